# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect refreshed data scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row => new value for column F
$sheetExhibitions = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    4  = 92
    5  = 52
    6  = 33
    7  = 582
    8  = 75
    9  = 8479
    10 = 790
    12 = 1129
    13 = 908
    14 = 79
    16 = 221
    17 = 186
    19 = 226
    20 = 956
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibitions.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (All types) - row => new value for column F
$sheetAllTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    4  = 92
    6  = 52
    7  = 33
    9  = 582
    10 = 75
    11 = 8479
    12 = 790
    14 = 1129
    15 = 908
    16 = 79
    18 = 221
    19 = 186
    21 = 226
    22 = 956
}
foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
